$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Strategies")

# Rename the strategy in A4 from "wma20_system" to "wma20_macd_system"
$ws.Range("A4").Value = "wma20_macd_system"

# Move the active selection to F16 (matches final cursor position on save)
$ws.Activate()
$ws.Range("F16").Select()
